# Adds a new "area" column (D) to the County/City Info sheet, shifting the
# existing "region" (old D) and "special_municipality" (old E) columns one
# position to the right (E and F respectively), and fills in area values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 23

# ---------------------------------------------------------------------
# Step 1: shift existing columns out of the way before overwriting them.
# Old layout: A=city_name B=city_code C=population D=region E=special_municipality
# New layout: A=city_name B=city_code C=population D=area E=region F=special_municipality
# Copy E -> F first (F is currently empty, so nothing is lost), then D -> E.
# Copy formats then values separately, which round-trips styles reliably.
# ---------------------------------------------------------------------

$rngE = $ws.Range("E1:E$lastRow")
$rngF = $ws.Range("F1:F$lastRow")
$rngE.Copy()
$rngF.PasteSpecial(-4122)   # xlPasteFormats
$rngE.Copy()
$rngF.PasteSpecial(-4163)   # xlPasteValues

$rngD = $ws.Range("D1:D$lastRow")
$rngD.Copy()
$rngE.PasteSpecial(-4122)   # xlPasteFormats
$rngD.Copy()
$rngE.PasteSpecial(-4163)   # xlPasteValues

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Step 2: set the new column D header to "area" (keep the existing bold
# header style that was already on D1).
# ---------------------------------------------------------------------

$ws.Range("D1").Value2 = "area"

# ---------------------------------------------------------------------
# Step 3: fill in the area values (sq. km) for each county/city, formatting
# each cell with the Arial / theme-1-color font used elsewhere in the sheet.
# Values under 1000 keep General formatting; values >= 1000 get a
# thousands-separated "#,##0.00" format (mirrors Excel's automatic
# formatting behavior for typed numbers).
# ---------------------------------------------------------------------

$areaByRow = @{
    2  = 271.79
    3  = 2052.57
    4  = 1220.95
    5  = 2214.9
    6  = 2191.65
    7  = 2951.85
    8  = 132.76
    9  = 104.15
    10 = 1427.59
    11 = 1820.31
    12 = 1074.39
    13 = 4106.43
    14 = 1290.84
    15 = 60.03
    16 = 1901.67
    17 = 2775.6
    18 = 2143.62
    19 = 4628.57
    20 = 3515.25
    21 = 141.05
    22 = 150.46
    23 = 29.6
}

foreach ($row in 2..$lastRow) {
    $value = $areaByRow[$row]
    $cell = $ws.Range("D$row")
    $cell.Font.Name = "Arial"
    $cell.Font.ThemeColor = 1
    if ($value -ge 1000) {
        $cell.NumberFormat = "#,##0.00"
    }
    $cell.Value2 = $value
}

Write-Output "area column inserted"
